# Weekly Fruta/Hortaliza price update: insert two new rows of data (Fortuna
# variety) at the top of the dated block (row 75) and push the existing
# rows down by two, extending the used range from T101 to T103.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above current row 75; this shifts rows 75:101 down
# to 77:103 and (as Excel normally does) copies the formatting - including
# the date number format in column D - from the row above into the new rows.
$ws.Rows.Item(75).Insert()
$ws.Rows.Item(75).Insert()

# New row 75: Ciruela / Fortuna / Primera
$ws.Cells.Item(75, 1).Value  = 11
$ws.Cells.Item(75, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(75, 3).Value  = "Bíobío"
$ws.Cells.Item(75, 4).Value  = 44964
$ws.Cells.Item(75, 5).Value  = 8
$ws.Cells.Item(75, 6).Value  = "Fruta"
$ws.Cells.Item(75, 7).Value  = 100103
$ws.Cells.Item(75, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(75, 9).Value  = 100103002
$ws.Cells.Item(75, 10).Value = "Ciruela"
$ws.Cells.Item(75, 11).Value = "Fortuna"
$ws.Cells.Item(75, 12).Value = "Primera"
$ws.Cells.Item(75, 13).Value = 100
$ws.Cells.Item(75, 14).Value = 11000
$ws.Cells.Item(75, 15).Value = 12000
$ws.Cells.Item(75, 16).Value = 11500
$ws.Cells.Item(75, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(75, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(75, 19).Value = 639
$ws.Cells.Item(75, 20).Value = 18

# New row 76: Ciruela / Fortuna / Segunda
$ws.Cells.Item(76, 1).Value  = 11
$ws.Cells.Item(76, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(76, 3).Value  = "Bíobío"
$ws.Cells.Item(76, 4).Value  = 44964
$ws.Cells.Item(76, 5).Value  = 8
$ws.Cells.Item(76, 6).Value  = "Fruta"
$ws.Cells.Item(76, 7).Value  = 100103
$ws.Cells.Item(76, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(76, 9).Value  = 100103002
$ws.Cells.Item(76, 10).Value = "Ciruela"
$ws.Cells.Item(76, 11).Value = "Fortuna"
$ws.Cells.Item(76, 12).Value = "Segunda"
$ws.Cells.Item(76, 13).Value = 50
$ws.Cells.Item(76, 14).Value = 9000
$ws.Cells.Item(76, 15).Value = 9000
$ws.Cells.Item(76, 16).Value = 9000
$ws.Cells.Item(76, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(76, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(76, 19).Value = 500
$ws.Cells.Item(76, 20).Value = 18
